$wb = $excel.ActiveWorkbook

# --- TextFileSequence sheet ---
# Before: A=sequence_file_format, B=file_content, C=overhang_crick_3prime, D=overhang_watson_3prime, E=id, F=type
# After:  A=sequence_file_format, B=overhang_crick_3prime, C=overhang_watson_3prime, D=file_content, E=id, F=type
$wsText = $wb.Worksheets.Item("TextFileSequence")
$wsText.Range("A1").Value = "sequence_file_format"
$wsText.Range("B1").Value = "overhang_crick_3prime"
$wsText.Range("C1").Value = "overhang_watson_3prime"
$wsText.Range("D1").Value = "file_content"
$wsText.Range("E1").Value = "id"
$wsText.Range("F1").Value = "type"

# --- ManuallyTypedSource sheet ---
# Before: A=user_input, B=circular, C=input, D=output, E=type, F=id
# After:  A=overhang_crick_3prime, B=overhang_watson_3prime, C=user_input, D=circular, E=input, F=output, G=type, H=id
$wsManual = $wb.Worksheets.Item("ManuallyTypedSource")
$wsManual.Range("A1").Value = "overhang_crick_3prime"
$wsManual.Range("B1").Value = "overhang_watson_3prime"
$wsManual.Range("C1").Value = "user_input"
$wsManual.Range("D1").Value = "circular"
$wsManual.Range("E1").Value = "input"
$wsManual.Range("F1").Value = "output"
$wsManual.Range("G1").Value = "type"
$wsManual.Range("H1").Value = "id"

# --- OligoHybridizationSource sheet ---
# Before: A=forward_oligo, B=reverse_oligo, C=overhang_crick_3prime, D=input, E=output, F=type, G=id
# After:  A=overhang_crick_3prime, B=forward_oligo, C=reverse_oligo, D=input, E=output, F=type, G=id
$wsOligo = $wb.Worksheets.Item("OligoHybridizationSource")
$wsOligo.Range("A1").Value = "overhang_crick_3prime"
$wsOligo.Range("B1").Value = "forward_oligo"
$wsOligo.Range("C1").Value = "reverse_oligo"
$wsOligo.Range("D1").Value = "input"
$wsOligo.Range("E1").Value = "output"
$wsOligo.Range("F1").Value = "type"
$wsOligo.Range("G1").Value = "id"
